# OO-4969: don't ask length of text entry, it's not used in excel import
#
# The "FIB" (fill-in-the-blank) help rows asked authors for both a gap
# length and a character-count ("Laenge,AnzahlZeichen", e.g. "20,50").
# The character-count part is unused by the importer, so the column now
# only wants a plain numeric gap length, and the help text in column D is
# updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Row 28 (first FIB example, "Eine Luecke ...") : was text "20,50" -> now
# a plain number (just the gap length, 20).
$ws.Range("C28").Value = 20

# Updated help text explaining the (now single) numeric value.
$ws.Range("D28").Value = "Eine Lücke, Länge 20. Wenn richtig gibt’s einen Punkt. Strickpunkt trennt Synonyme bzw. mehrere korrekte Lösungen"

# Rows 35, 37 and 39 (second FIB example, multiple gaps) : each gap's
# "Laenge,AnzahlZeichen" text "2,2" becomes the plain number 2.
$ws.Range("C35").Value = 2
$ws.Range("C37").Value = 2
$ws.Range("C39").Value = 2

# Match the author's final cursor position/selection on the sheet.
$ws.Activate()
$ws.Range("D28").Select()
